$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 16; $r++) {
    $wsOverview.Range("G$r").Value = "2016-08-20 18:35:49"
}

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
for ($r = 2; $r -le 16; $r++) {
    $wsZhCn.Range("E$r").Value = "mt"
    $wsZhCn.Range("H$r").Value = "2016-08-20 18:35:44"
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
for ($r = 2; $r -le 16; $r++) {
    $wsDeDe.Range("E$r").Value = "mt"
    $wsDeDe.Range("H$r").Value = "2016-08-20 18:35:49"
}
